$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Update counts for the 50 scale FOCUS models (Comments column) ---
$ws.Range("I21").Value = "24 complete"
$ws.Range("I23").Value = "12 complete"

# --- Add new ACHD model print requests (rows 31-33) ---
$ws.Range("A31").Value = "30-05-2018"
$ws.Range("C31").Value = "Full Fontan"
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = "Polylite"
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 20
$ws.Range("H31").Value = 0.2
$ws.Range("I31").Value = "NA"

$ws.Range("A32").Value = "30-05-2018"
$ws.Range("C32").Value = "Full Dextrocardia"
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = "Polylite"
$ws.Range("F32").Value = 2
$ws.Range("G32").Value = 20
$ws.Range("H32").Value = 0.2
$ws.Range("I32").Value = "NA"

$ws.Range("A33").Value = "30-05-2018"
$ws.Range("C33").Value = "Full TOF"
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = "Polylite"
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = 20
$ws.Range("H33").Value = 0.2
$ws.Range("I33").Value = "NA"

# --- Best-effort view-state: leave the cursor on the last touched row ---
$ws.Range("I25").Select()
